# act tablas web jul25
# Updates the "220103" indicator table: adds 2023/2022 at the top of the
# Data sheet, extends the historical series back from 2004 to 1985 (1989
# is missing from the source data) at the bottom, refreshes several of the
# existing "Valor" figures, and records the "actualizacion" metadata row.

function Set-TextValue($cell, [string]$text) {
    # Writing a digit-only string through .Value lets the host coerce it to
    # a number. Forcing a Text number format keeps it a string, then
    # clearing formats afterwards drops the leftover style so the cell ends
    # up with default formatting again (matching the source workbook).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Data sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Data")

# Insert two fresh rows under the header for the new 2023 / 2022 figures;
# this pushes the existing 2021..2005 rows down by two without touching
# their values.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

Set-TextValue $ws.Cells.Item(2, 1) "2023"
$ws.Cells.Item(2, 2).Value = 5.2

Set-TextValue $ws.Cells.Item(3, 1) "2022"
$ws.Cells.Item(3, 2).Value = 5.4

# Extend the series backwards from 2004 to 1985 (note: 1989 is absent from
# the source), appended after the row that now holds 2005 (row 20).
$history = @(
    @("2004", 7.1),
    @("2003", 6.8),
    @("2002", 6.4),
    @("2001", 7.1),
    @("2000", 5.8),
    @("1999", 6),
    @("1998", 5.7),
    @("1997", 6.2),
    @("1996", 6.3),
    @("1995", 3.8),
    @("1994", 4.1),
    @("1993", 3.8),
    @("1992", 3.6),
    @("1991", 2.3),
    @("1990", 2.9),
    @("1988", 3.7),
    @("1987", 4.8),
    @("1986", 4.7),
    @("1985", 4.9)
)

$row = 21
foreach ($entry in $history) {
    Set-TextValue $ws.Cells.Item($row, 1) $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Metadata")

# A1 moves from an empty string to a single space (matching B1).
$ws2.Cells.Item(1, 1).Value = " "

# Insert the "actualizacion" / "Julio 2025" row just above "cita".
$ws2.Rows.Item(9).Insert()
$ws2.Cells.Item(9, 1).Value = "actualizacion"
$ws2.Cells.Item(9, 2).Value = "Julio 2025"
